# FIX TIMEDATE & CREATE ROOM(block same name)
#
# Adds SUCCESS/FAIL/NAME_FAIL status-code + Korean message pairs (columns I/J)
# to the three "room"-related API blocks (rows 43-44, 47-48, 51-53), widens
# column J slightly, and updates the saved view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- group_create_check.php block (row 43-44) ---------------------------
$ws.Range("I43").Value = "SUCCESS"
$ws.Range("J43").Value = "성공적으로 생성됨 알림"
$ws.Range("I44").Value = "FAIL"
$ws.Range("J44").Value = "실패 알림"

# --- second block (row 47-48) --------------------------------------------
$ws.Range("I47").Value = "SUCCESS"
$ws.Range("J47").Value = "성공적으로 생성됨 알림"
$ws.Range("I48").Value = "FAIL"
$ws.Range("J48").Value = "실패 알림"

# --- third block (row 51-53) - also guards against duplicate room names --
$ws.Range("I51").Value = "SUCCESS"
$ws.Range("J51").Value = "성공적으로 생성됨 알림"
$ws.Range("I52").Value = "FAIL"
$ws.Range("J52").Value = "실패 알림"
$ws.Range("I53").Value = "NAME_FAIL"
$ws.Range("J53").Value = "이미 해당 방 이름이 존재"

# --- widen column J (now holds the longer Korean notification text) -----
$ws.Range("J1").ColumnWidth = 20.75

# --- move the saved viewport / selection down to the new rows -----------
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A58").Select()
